$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "309.35"
Set-TextValue $ws "E2" "-0.48%"
Set-TextValue $ws "G2" "14"

Set-TextValue $ws "D3" "37.34"
Set-TextValue $ws "E3" "-0.72%"
Set-TextValue $ws "G3" "14"

Set-TextValue $ws "D4" "5.128"
Set-TextValue $ws "E4" "0.27%"
Set-TextValue $ws "G4" "14"

Set-TextValue $ws "D5" "0.07847"
Set-TextValue $ws "E5" "-0.58%"
Set-TextValue $ws "G5" "14"

Set-TextValue $ws "B6" "KuCoinToken"
Set-TextValue $ws "C6" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws "D6" "8.265"
Set-TextValue $ws "E6" "0.47%"
Set-TextValue $ws "G6" "14"

Set-TextValue $ws "B7" "FTXToken"
Set-TextValue $ws "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D7" "1.888"
Set-TextValue $ws "E7" "-0.82%"
Set-TextValue $ws "G7" "14"

Set-TextValue $ws "B8" "BTSEToken"
Set-TextValue $ws "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D8" "2.943"
Set-TextValue $ws "E8" "2.76%"
Set-TextValue $ws "G8" "14"

Set-TextValue $ws "B9" "MXToken"
Set-TextValue $ws "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D9" "0.9239"
Set-TextValue $ws "E9" "-0.26%"
Set-TextValue $ws "G9" "14"

Set-TextValue $ws "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D10" "0.1146"
Set-TextValue $ws "E10" "-4.46%"
Set-TextValue $ws "G10" "14"

Set-TextValue $ws "B11" "WazirX"
Set-TextValue $ws "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D11" "0.1906"
Set-TextValue $ws "E11" "-0.17%"
Set-TextValue $ws "G11" "14"

Set-TextValue $ws "B12" "MandalaExchangeToken"
Set-TextValue $ws "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D12" "0.08846"
Set-TextValue $ws "E12" "-6.65%"
Set-TextValue $ws "G12" "14"

Set-TextValue $ws "B13" "BitrueCoin"
Set-TextValue $ws "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D13" "0.03333"
Set-TextValue $ws "E13" "-1.00%"
Set-TextValue $ws "G13" "14"

Set-TextValue $ws "B14" "BitMartToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D14" "0.09618"
Set-TextValue $ws "E14" "-0.07%"
Set-TextValue $ws "G14" "14"

Set-TextValue $ws "B15" "BitForexToken"
Set-TextValue $ws "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D15" "0.001385"
Set-TextValue $ws "E15" "1.29%"
Set-TextValue $ws "G15" "14"

Set-TextValue $ws "B16" "TigerCash"
Set-TextValue $ws "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D16" "0.006002"
Set-TextValue $ws "E16" "2.06%"
Set-TextValue $ws "G16" "14"

Set-TextValue $ws "B17" "LEO"
Set-TextValue $ws "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D17" "3.392"
Set-TextValue $ws "E17" "-3.90%"
Set-TextValue $ws "G17" "14"

Set-TextValue $ws "B18" "GateToken"
Set-TextValue $ws "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D18" "4.392"
Set-TextValue $ws "E18" "-0.19%"
Set-TextValue $ws "G18" "14"

Set-TextValue $ws "E19" "0.78%"
Set-TextValue $ws "G19" "14"

Set-TextValue $ws "D20" "6.330"
Set-TextValue $ws "E20" "20.42%"
Set-TextValue $ws "G20" "14"

Set-TextValue $ws "D21" "0.1297"
Set-TextValue $ws "E21" "1.81%"
Set-TextValue $ws "G21" "14"

Set-TextValue $ws "E22" "-6.95%"
Set-TextValue $ws "G22" "14"

Set-TextValue $ws "D23" "0.04355"
Set-TextValue $ws "E23" "-0.09%"
Set-TextValue $ws "G23" "14"

Set-TextValue $ws "E24" "0.66%"
Set-TextValue $ws "G24" "14"

Set-TextValue $ws "D25" "0.004275"
Set-TextValue $ws "E25" "0.00%"
Set-TextValue $ws "G25" "14"

Set-TextValue $ws "D26" "0.0001401"
Set-TextValue $ws "E26" "8.04%"
Set-TextValue $ws "G26" "14"

Set-TextValue $ws "D27" "0.0002902"
Set-TextValue $ws "G27" "14"

Set-TextValue $ws "G28" "14"

Set-TextValue $ws "G29" "14"

Set-TextValue $ws "G30" "14"

Set-TextValue $ws "G31" "14"

Set-TextValue $ws "G32" "14"

Set-TextValue $ws "G33" "14"

Set-TextValue $ws "G34" "14"

Set-TextValue $ws "G35" "14"

Set-TextValue $ws "G36" "14"

Set-TextValue $ws "G37" "14"

Set-TextValue $ws "G38" "14"

Set-TextValue $ws "D39" "0.02166"
Set-TextValue $ws "E39" "3.42%"
Set-TextValue $ws "G39" "14"

Set-TextValue $ws "D40" "0.05012"
Set-TextValue $ws "E40" "-1.82%"
Set-TextValue $ws "G40" "14"

Set-TextValue $ws "D41" "0.007580"
Set-TextValue $ws "E41" "-0.27%"
Set-TextValue $ws "G41" "14"

Set-TextValue $ws "D42" "0.1354"
Set-TextValue $ws "E42" "0.15%"
Set-TextValue $ws "G42" "14"

Set-TextValue $ws "D43" "0.008483"
Set-TextValue $ws "E43" "-6.83%"
Set-TextValue $ws "G43" "14"

Set-TextValue $ws "D44" "0.002008"
Set-TextValue $ws "E44" "-2.30%"
Set-TextValue $ws "G44" "14"

Set-TextValue $ws "E45" "-5.30%"
Set-TextValue $ws "G45" "14"

Set-TextValue $ws "D46" "0.00006570"
Set-TextValue $ws "E46" "-1.70%"
Set-TextValue $ws "G46" "14"

Set-TextValue $ws "E47" "0.38%"
Set-TextValue $ws "G47" "14"

Set-TextValue $ws "D48" "0.003295"
Set-TextValue $ws "E48" "14.32%"
Set-TextValue $ws "G48" "14"

Set-TextValue $ws "D49" "0.001444"
Set-TextValue $ws "E49" "20.71%"
Set-TextValue $ws "G49" "14"

Set-TextValue $ws "D50" "0.00002101"
Set-TextValue $ws "E50" "0.38%"
Set-TextValue $ws "G50" "14"

Set-TextValue $ws "D51" "0.0002001"
Set-TextValue $ws "E51" "0.38%"
Set-TextValue $ws "G51" "14"
